# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Ajo" (Terminal La Palmera de La Serena)
# as a new row 148, shifting the existing rows 148:196 down to 149:197.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 148:196 down one row, opening up a blank row 148.
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new weekly record.
$ws.Range("A148").Value = 8
$ws.Range("B148").Value = "Terminal La Palmera de La Serena"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44559
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 100112003
$ws.Range("G148").Value = "Ajo"
$ws.Range("H148").Value = "Chino"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 560
$ws.Range("K148").Value = 19000
$ws.Range("L148").Value = 20000
$ws.Range("M148").Value = 19500
$ws.Range("N148").Value = "$/caja 10 kilos"
$ws.Range("O148").Value = "China"
$ws.Range("P148").Value = 1950
$ws.Range("Q148").Value = 10
$ws.Range("R148").Value = "Hortaliza"
